$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New B, C, D, E, G values per row (G = B + C + D + E ; F "Win" column unchanged)
$data = @{
    2  = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694)
    3  = @(3.286832544864788, 0.306821227259698, 3.537761648806719, 0.4942365360607697, 7.625651956991975)
    4  = @(0.04271373187048222, 0.04071648406533734, 3.537761648806719, 0.4942365360607697, 4.115428400803308)
    5  = @(0.6606524410359556, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 6.348428708163715)
    6  = @(1.455362044514542, 0.04071648406533734, 0.1494219747398047, 0.4942365360607697, 2.139737039380454)
    7  = @(1.455362044514542, 1.655778082260271, 0.7527432677738641, 10.19245300693656, 14.05633640148523)
    8  = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 5.586269137925634)
    9  = @(0.1190320826869504, 0.04071648406533734, 3.537761648806719, 0.4942365360607697, 4.191746751619776)
    10 = @(0.6606524410359556, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 2.960089034096801)
    11 = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 5.586269137925634)
    12 = @(1.455362044514542, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 3.754798637575387)
    13 = @(0.1190320826869504, 0.306821227259698, 0.7527432677738641, 0.4942365360607697, 1.672833113781282)
    14 = @(1.455362044514542, 1.655778082260271, 22.3905356188092, 0.4942365360607697, 25.99591228164478)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 2).Value = $vals[0]
    $ws.Cells.Item($r, 3).Value = $vals[1]
    $ws.Cells.Item($r, 4).Value = $vals[2]
    $ws.Cells.Item($r, 5).Value = $vals[3]
    $ws.Cells.Item($r, 7).Value = $vals[4]
}
